function Wrap-Xml($inner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><w:body>' + $inner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

# =========================================================================
# 1) "файл doxc" figure paragraph (currently paragraph 23): split the
#    BodyText paragraph (image + "{#fig:004 width=70%}" text) into a
#    CaptionedFigure paragraph (image only, title="fig:") followed by an
#    ImageCaption paragraph with the figure's caption text.
# =========================================================================
$p = $d.Paragraphs(23)
$inner = '<w:p><w:pPr><w:pStyle w:val="CaptionedFigure"/></w:pPr><w:r><w:drawing><wp:inline><wp:extent cx="3733800" cy="86786"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="файл doxc" title="fig:" id="33" name="Picture"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="image/4.png" id="34" name="Picture"/><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId32"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="3733800" cy="86786"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ImageCaption"/></w:pPr><w:r><w:t xml:space="preserve">файл doxc</w:t></w:r></w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# =========================================================================
# 2) "файл pdf" figure paragraph (now paragraph 26, since step 1 added one
#    extra paragraph): same split as above.
# =========================================================================
$p = $d.Paragraphs(26)
$inner = '<w:p><w:pPr><w:pStyle w:val="CaptionedFigure"/></w:pPr><w:r><w:drawing><wp:inline><wp:extent cx="3733800" cy="86240"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="файл pdf" title="fig:" id="36" name="Picture"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="image/5.png" id="37" name="Picture"/><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId35"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="3733800" cy="86240"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ImageCaption"/></w:pPr><w:r><w:t xml:space="preserve">файл pdf</w:t></w:r></w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# =========================================================================
# 3) Fix the stray markdown-footnote-ish typo: "(рис. ??]." -> "(рис. -??)."
# =========================================================================
$d.Content.Find.Execute("doxc файл (рис. ??].", $true, $false, $false, $false, $false, $true, 1, $false, "doxc файл (рис. -??).", 2) | Out-Null

Write-Output ("paragraphs after part 1-3=" + $d.Paragraphs.Count)
for ($i = 22; $i -le 29; $i++) {
    $pp = $d.Paragraphs($i)
    Write-Output ($i.ToString() + "|" + $pp.Style.NameLocal + "|[" + $pp.Range.Text + "]")
}
